# Apply the data edits described by the commit:
# "Retry failed testcases and IAnnotation Transformer to remove annotations from @Test"
#
# RUNMANAGER sheet: stop executing "loginlogoutTest" (C2: yes -> no)
# DATA sheet: enable execution for the "newtest" / admin row (B5: no -> yes)

$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")

# RUNMANAGER!C2 : yes -> no
$wsRunManager.Range("C2").Value = "no"

# DATA!B5 : no -> yes
$wsData.Range("B5").Value = "yes"

# Move the active selection on the DATA sheet to C5, matching the cell that
# was just edited, and make DATA the active sheet.
$wsData.Activate()
$wsData.Range("C5").Select()
